# documentacion.docx - "implementando estilos globales y por componente"
#
# 1) The original document has a number of Word auto spell/grammar-check
#    markers (<w:proofErr .../>) splitting otherwise-contiguous text into
#    several <w:r> runs (e.g. "npm" / " i @angular/cli", "server" /
#    "-side" / " rendering", etc). Re-reading a range's WordOpenXML and
#    feeding it straight back in via InsertXML makes Word re-flow that
#    range from scratch, which drops the now-stale proofing marks and
#    coalesces the text back into minimal runs - exactly the cleanup seen
#    in the diff - without touching the actual wording, the table, or the
#    comment anchored on "--skip-tests".
$d = $word.ActiveDocument

$full = $d.Content
$full.InsertXML($full.WordOpenXML)

# 2) The trailing empty paragraph after "Uso del router-outlet" gets the
#    gist link as its text.
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Range.Text = "https://gist.github.com/nicobytes/ba2252b0b5ac2cbdafc40c0accd24862"
